# "Update latest lessons file"
# Rename the worksheet and move the saved selection/active cell,
# matching the author's re-purposing of this workbook as the
# canonical "Lessons" file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet tab/name: "perfect can u generate for 3 mo" -> "Lessons"
$ws.Name = "Lessons"

# Saved cursor position: G168:G169 -> B6 (single cell)
$ws.Activate()
$ws.Range("B6").Select()
